# Fix the pptx logo position: move the logo picture from the right side
# of the slide back to the left side (x offset 7863840 EMU -> 182880 EMU).
# 182880 EMU = 0.2 in = 14.4 pt (matches the y-offset already used).

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.Name -eq "Picture 3" -and [math]::Round($sh.Left, 2) -eq 619.2) {
            # 182880 EMU = 14.4 pt. Use a value that round-trips precisely
            # through the COM single-precision float marshalling so the
            # saved XML offset is exactly 182880 (not 182879).
            $sh.Left = 14.400024
        }
    }
}
